$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("observations")

$ws.Cells.Item(338, 1).Value = 'great horned owl'
$ws.Cells.Item(338, 2).Value = 'birds'
$ws.Cells.Item(338, 3).Value = 43639
$ws.Cells.Item(338, 4).Value = 0.75
$ws.Cells.Item(338, 6).Value = 'Maple Grove'
$ws.Cells.Item(338, 11).Value = 'approximate: "12 hours after thunderstorm"'

$ws.Cells.Item(339, 1).Value = 'american toad'
$ws.Cells.Item(339, 2).Value = 'herps'
$ws.Cells.Item(339, 3).Value = 43631
$ws.Cells.Item(339, 4).Value = 0.3333333333333333
$ws.Cells.Item(339, 6).Value = 'Otis'
$ws.Cells.Item(339, 11).Value = 'guessing at time: walking with brooklyn on our walk around the block, which I assume was our morning walk.'

$ws.Cells.Item(340, 1).Value = 'red-eyed vireo'
$ws.Cells.Item(340, 2).Value = 'birds'
$ws.Cells.Item(340, 3).Value = 43633
$ws.Cells.Item(340, 4).Value = 0.2708333333333333
$ws.Cells.Item(340, 6).Value = 'Maple Grove'
$ws.Cells.Item(340, 11).Value = 'guessing at time: morning, and it would have been b/f work'

$ws.Cells.Item(341, 1).Value = 'blue-gray gnatcatcher'
$ws.Cells.Item(341, 2).Value = 'birds'
$ws.Cells.Item(341, 3).Value = 43633
$ws.Cells.Item(341, 4).Value = 0.2708333333333333
$ws.Cells.Item(341, 6).Value = 'Maple Grove'
$ws.Cells.Item(341, 11).Value = 'guessing at time: morning, and it would have been b/f work'

$ws.Cells.Item(342, 1).Value = 'eastern wood peewee'
$ws.Cells.Item(342, 2).Value = 'birds'
$ws.Cells.Item(342, 3).Value = 43633
$ws.Cells.Item(342, 4).Value = 0.2708333333333333
$ws.Cells.Item(342, 6).Value = 'Maple Grove'
$ws.Cells.Item(342, 11).Value = 'guessing at time: morning, and it would have been b/f work'

$ws.Cells.Item(343, 1).Value = 'red-bellied woodpecker'
$ws.Cells.Item(343, 2).Value = 'birds'
$ws.Cells.Item(343, 3).Value = 43633
$ws.Cells.Item(343, 4).Value = 0.2708333333333333
$ws.Cells.Item(343, 6).Value = 'Maple Grove'
$ws.Cells.Item(343, 11).Value = 'guessing at time: morning, and it would have been b/f work'

$ws.Cells.Item(344, 1).Value = 'great crested flycatcher'
$ws.Cells.Item(344, 2).Value = 'birds'
$ws.Cells.Item(344, 3).Value = 43633
$ws.Cells.Item(344, 4).Value = 0.2708333333333333
$ws.Cells.Item(344, 6).Value = 'Maple Grove'
$ws.Cells.Item(344, 11).Value = 'guessing at time: morning, and it would have been b/f work'

$ws.Cells.Item(345, 1).Value = 'great crested flycatcher'
$ws.Cells.Item(345, 2).Value = 'birds'
$ws.Cells.Item(345, 3).Value = 43623
$ws.Cells.Item(345, 4).Value = 0.2916666666666667
$ws.Cells.Item(345, 6).Value = 'Morton Arboretum'
$ws.Cells.Item(345, 11).Value = 'guessing at time: morning ride into work, and I recall the morning, and I don''t think I was particularly early into work'

$ws.Cells.Item(346, 1).Value = 'eastern wood peewee'
$ws.Cells.Item(346, 2).Value = 'birds'
$ws.Cells.Item(346, 3).Value = 43623
$ws.Cells.Item(346, 4).Value = 0.2916666666666667
$ws.Cells.Item(346, 6).Value = 'Morton Arboretum'
$ws.Cells.Item(346, 11).Value = 'guessing at time: morning ride into work, and I recall the morning, and I don''t think I was particularly early into work'

$ws.Cells.Item(347, 1).Value = 'red-eyed vireo'
$ws.Cells.Item(347, 2).Value = 'birds'
$ws.Cells.Item(347, 3).Value = 43623
$ws.Cells.Item(347, 4).Value = 0.2916666666666667
$ws.Cells.Item(347, 6).Value = 'Morton Arboretum'
$ws.Cells.Item(347, 11).Value = 'guessing at time: morning ride into work, and I recall the morning, and I don''t think I was particularly early into work'

$ws.Cells.Item(348, 1).Value = 'great horned owl'
$ws.Cells.Item(348, 2).Value = 'birds'
$ws.Cells.Item(348, 3).Value = 43630
$ws.Cells.Item(348, 4).Value = 0.75
$ws.Cells.Item(348, 6).Value = 'Maple Grove'
$ws.Cells.Item(348, 11).Value = 'guessing at time: evening walk with Rachel, and it was light enough to see Carex davisii'

$ws.Cells.Item(349, 1).Value = 'cicada'
$ws.Cells.Item(349, 2).Value = 'insects'
$ws.Cells.Item(349, 3).Value = 43647
$ws.Cells.Item(349, 4).Value = 0.6666666666666666
$ws.Cells.Item(349, 6).Value = 'Otis'
$ws.Cells.Item(349, 11).Value = 'guessing at time: "late afternoon"'

$ws.Cells.Item(350, 1).Value = 'eastern wood peewee'
$ws.Cells.Item(350, 2).Value = 'birds'
$ws.Cells.Item(350, 3).Value = 43648
$ws.Cells.Item(350, 4).Value = 0.22916666666666666
$ws.Cells.Item(350, 6).Value = 'Morton Arboretum'
$ws.Cells.Item(350, 11).Value = '"just after sunrise"'

$ws.Cells.Item(351, 1).Value = 'indigo bunting'
$ws.Cells.Item(351, 2).Value = 'birds'
$ws.Cells.Item(351, 3).Value = 43648
$ws.Cells.Item(351, 4).Value = 0.22916666666666666
$ws.Cells.Item(351, 6).Value = 'Morton Arboretum'
$ws.Cells.Item(351, 11).Value = '"just after sunrise"'

$ws.Cells.Item(352, 1).Value = 'robin'
$ws.Cells.Item(352, 2).Value = 'birds'
$ws.Cells.Item(352, 3).Value = 43648
$ws.Cells.Item(352, 4).Value = 0.22916666666666666
$ws.Cells.Item(352, 6).Value = 'Morton Arboretum'
$ws.Cells.Item(352, 11).Value = '"just after sunrise"'

$ws.Cells.Item(353, 1).Value = 'house wren'
$ws.Cells.Item(353, 2).Value = 'birds'
$ws.Cells.Item(353, 3).Value = 43648
$ws.Cells.Item(353, 4).Value = 0.22916666666666666
$ws.Cells.Item(353, 6).Value = 'Morton Arboretum'
$ws.Cells.Item(353, 11).Value = '"just after sunrise"'

$ws.Cells.Item(354, 1).Value = 'red-eyed vireo'
$ws.Cells.Item(354, 2).Value = 'birds'
$ws.Cells.Item(354, 3).Value = 43648
$ws.Cells.Item(354, 4).Value = 0.22916666666666666
$ws.Cells.Item(354, 6).Value = 'Morton Arboretum'
$ws.Cells.Item(354, 11).Value = '"just after sunrise"'

$ws.Cells.Item(355, 1).Value = 'white-breasted nuthatch'
$ws.Cells.Item(355, 2).Value = 'birds'
$ws.Cells.Item(355, 3).Value = 43648
$ws.Cells.Item(355, 4).Value = 0.22916666666666666
$ws.Cells.Item(355, 6).Value = 'Morton Arboretum'
$ws.Cells.Item(355, 11).Value = '"just after sunrise"'

$ws.Cells.Item(356, 1).Value = 'field sparrow'
$ws.Cells.Item(356, 2).Value = 'birds'
$ws.Cells.Item(356, 3).Value = 43648
$ws.Cells.Item(356, 4).Value = 0.22916666666666666
$ws.Cells.Item(356, 6).Value = 'Morton Arboretum'
$ws.Cells.Item(356, 11).Value = '"just after sunrise"'

$ws.Cells.Item(357, 1).Value = 'white-breasted nuthatch'
$ws.Cells.Item(357, 2).Value = 'birds'
$ws.Cells.Item(357, 3).Value = 43697
$ws.Cells.Item(357, 4).Value = 0.8125
$ws.Cells.Item(357, 6).Value = 'Maple Grove'

$ws.Cells.Item(358, 1).Value = 'eastern wood peewee'
$ws.Cells.Item(358, 2).Value = 'birds'
$ws.Cells.Item(358, 3).Value = 43697
$ws.Cells.Item(358, 4).Value = 0.8229166666666666
$ws.Cells.Item(358, 6).Value = 'Maple Grove'

$ws.Cells.Item(359, 1).Value = 'crickets'
$ws.Cells.Item(359, 2).Value = 'insects'
$ws.Cells.Item(359, 3).Value = 43697
$ws.Cells.Item(359, 4).Value = 0.8333333333333334
$ws.Cells.Item(359, 6).Value = 'Maple Grove'

$ws.Cells.Item(360, 1).Value = 'cicada'
$ws.Cells.Item(360, 2).Value = 'insects'
$ws.Cells.Item(360, 3).Value = 43697
$ws.Cells.Item(360, 4).Value = 0.8333333333333334
$ws.Cells.Item(360, 6).Value = 'Maple Grove'

$ws.Cells.Item(361, 1).Value = 'eastern wood peewee'
$ws.Cells.Item(361, 2).Value = 'birds'
$ws.Cells.Item(361, 3).Value = 43697
$ws.Cells.Item(361, 4).Value = 0.8333333333333334
$ws.Cells.Item(361, 6).Value = 'Maple Grove'

$ws.Cells.Item(362, 1).Value = 'crickets'
$ws.Cells.Item(362, 2).Value = 'insects'
$ws.Cells.Item(362, 3).Value = 43700
$ws.Cells.Item(362, 4).Value = 0.6875
$ws.Cells.Item(362, 6).Value = 'Maple Grove'

$ws.Cells.Item(363, 1).Value = 'eastern phoebe'
$ws.Cells.Item(363, 2).Value = 'birds'
$ws.Cells.Item(363, 3).Value = 43700
$ws.Cells.Item(363, 4).Value = 0.7083333333333334
$ws.Cells.Item(363, 6).Value = 'Maple Grove'

$ws.Cells.Item(364, 1).Value = 'song sparrow'
$ws.Cells.Item(364, 2).Value = 'birds'
$ws.Cells.Item(364, 3).Value = 43736
$ws.Cells.Item(364, 4).Value = 0.2708333333333333
$ws.Cells.Item(364, 6).Value = 'Mays Lake'

$ws.Cells.Item(365, 1).Value = 'eastern wood peewee'
$ws.Cells.Item(365, 2).Value = 'birds'
$ws.Cells.Item(365, 3).Value = 43721
$ws.Cells.Item(365, 4).Value = 0.2916666666666667
$ws.Cells.Item(365, 6).Value = 'Morton Arboretum'
$ws.Cells.Item(365, 11).Value = 'guessing at time: ride into work'

$ws.Cells.Item(366, 1).Value = 'crickets'
$ws.Cells.Item(366, 2).Value = 'insects'
$ws.Cells.Item(366, 3).Value = 43721
$ws.Cells.Item(366, 4).Value = 0.2916666666666667
$ws.Cells.Item(366, 6).Value = 'Morton Arboretum'
$ws.Cells.Item(366, 11).Value = 'guessing at time: ride into work'

$ws.Cells.Item(367, 1).Value = 'katydid'
$ws.Cells.Item(367, 2).Value = 'insects'
$ws.Cells.Item(367, 3).Value = 43722
$ws.Cells.Item(367, 4).Value = 0.8333333333333334
$ws.Cells.Item(367, 6).Value = 'Lillstreet art center'
$ws.Cells.Item(367, 11).Value = 'guessing at time: it was getting dark out, but I recall that we could still see a bit'

$ws.Cells.Item(368, 1).Value = 'woodcock'
$ws.Cells.Item(368, 2).Value = 'birds'
$ws.Cells.Item(368, 3).Value = 43755
$ws.Cells.Item(368, 4).Value = 0.8333333333333334
$ws.Cells.Item(368, 6).Value = 'Morton Arboretum'
$ws.Cells.Item(368, 11).Value = 'guessing at time: it was a report from Robb Telfer'

$ws.Cells.Item(369, 1).Value = 'crickets'
$ws.Cells.Item(369, 3).Value = 43738
$ws.Cells.Item(369, 4).Value = 0.2916666666666667
$ws.Cells.Item(369, 6).Value = 'Morton Arboretum'
$ws.Cells.Item(369, 11).Value = 'guessing at time: bike ride in'

$ws.Cells.Item(370, 1).Value = 'spring peeper'
$ws.Cells.Item(370, 2).Value = 'herps'
$ws.Cells.Item(370, 3).Value = 43738
$ws.Cells.Item(370, 4).Value = 0.2916666666666667
$ws.Cells.Item(370, 6).Value = 'Morton Arboretum'
$ws.Cells.Item(370, 11).Value = 'bracketing time -- "spring peepers were singing all day"'

$ws.Cells.Item(371, 1).Value = 'spring peeper'
$ws.Cells.Item(371, 2).Value = 'herps'
$ws.Cells.Item(371, 3).Value = 43738
$ws.Cells.Item(371, 4).Value = 0.6875
$ws.Cells.Item(371, 6).Value = 'Morton Arboretum'

$ws.Cells.Item(372, 1).Value = 'spring peeper'
$ws.Cells.Item(372, 2).Value = 'herps'
$ws.Cells.Item(372, 3).Value = 43740
$ws.Cells.Item(372, 4).Value = 0.3125
$ws.Cells.Item(372, 6).Value = 'Morton Arboretum'

$ws.Cells.Item(373, 1).Value = 'woodcock'
$ws.Cells.Item(373, 2).Value = 'birds'
$ws.Cells.Item(373, 3).Value = 43770
$ws.Cells.Item(373, 4).Value = 0.2916666666666667
$ws.Cells.Item(373, 6).Value = 'Morton Arboretum'

$ws.Cells.Item(374, 1).Value = 'junco'
$ws.Cells.Item(374, 2).Value = 'birds'
$ws.Cells.Item(374, 3).Value = 43770
$ws.Cells.Item(374, 4).Value = 0.2916666666666667
$ws.Cells.Item(374, 6).Value = 'Morton Arboretum'

$ws.Cells.Item(375, 1).Value = 'robin'
$ws.Cells.Item(375, 2).Value = 'birds'
$ws.Cells.Item(375, 3).Value = 43770
$ws.Cells.Item(375, 4).Value = 0.2916666666666667
$ws.Cells.Item(375, 6).Value = 'Morton Arboretum'

$ws.Cells.Item(376, 1).Value = 'great horned owl'
$ws.Cells.Item(376, 2).Value = 'birds'
$ws.Cells.Item(376, 3).Value = 43784
$ws.Cells.Item(376, 4).Value = 0.7083333333333334
$ws.Cells.Item(376, 6).Value = 'Maple Grove'

$ws.Cells.Item(377, 1).Value = 'sandhill crane'
$ws.Cells.Item(377, 2).Value = 'birds'
$ws.Cells.Item(377, 3).Value = 43784
$ws.Cells.Item(377, 4).Value = 0.625
$ws.Cells.Item(377, 6).Value = 'Otis'
$ws.Cells.Item(377, 11).Value = 'guessing at time: afternoon'

$ws.Cells.Item(378, 1).Value = 'white-breasted nuthatch'
$ws.Cells.Item(378, 2).Value = 'birds'
$ws.Cells.Item(378, 3).Value = 43784
$ws.Cells.Item(378, 4).Value = 0.6458333333333334
$ws.Cells.Item(378, 6).Value = 'Maple Grove'
$ws.Cells.Item(378, 11).Value = 'an hour b/f sunset'

$ws.Cells.Item(379, 1).Value = 'hairy woodpecker'
$ws.Cells.Item(379, 2).Value = 'birds'
$ws.Cells.Item(379, 3).Value = 43784
$ws.Cells.Item(379, 4).Value = 0.6458333333333334
$ws.Cells.Item(379, 6).Value = 'Maple Grove'

$ws.Cells.Item(380, 1).Value = 'red-bellied woodpecker'
$ws.Cells.Item(380, 2).Value = 'birds'
$ws.Cells.Item(380, 3).Value = 43784
$ws.Cells.Item(380, 4).Value = 0.6458333333333334
$ws.Cells.Item(380, 6).Value = 'Maple Grove'

$ws.Cells.Item(381, 1).Value = 'great horned owl'
$ws.Cells.Item(381, 2).Value = 'birds'
$ws.Cells.Item(381, 3).Value = 43784
$ws.Cells.Item(381, 4).Value = 0.6875
$ws.Cells.Item(381, 6).Value = 'Maple Grove'

$ws.Cells.Item(382, 1).Value = 'white-breasted nuthatch'
$ws.Cells.Item(382, 2).Value = 'birds'
$ws.Cells.Item(382, 3).Value = 43816
$ws.Cells.Item(382, 4).Value = 0.2916666666666667
$ws.Cells.Item(382, 6).Value = 'Maple Grove'
$ws.Cells.Item(382, 11).Value = 'guessing at time -- morning walk w/ dog on a Saturday'

$ws.Cells.Item(383, 1).Value = 'blue jay'
$ws.Cells.Item(383, 2).Value = 'birds'
$ws.Cells.Item(383, 3).Value = 43816
$ws.Cells.Item(383, 4).Value = 0.2916666666666667
$ws.Cells.Item(383, 6).Value = 'Maple Grove'
$ws.Cells.Item(383, 11).Value = 'guessing at time -- morning walk w/ dog on a Saturday'
$ws.Range("H382").Select()
